# "Generate Report for Handback" — mark a.md / b.md as handed back and
# in sync with en-US, and record the handback report (target file,
# handback file, handback datetime) for the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# Cornflower blue (FF6495ED), matching the workbook's existing HyperLink style.
$hyperlinkColor = 15570276

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Overview sheet: both locale status columns move from "Ready for
# handoff" to "Handed back: in sync with en-US" for the two real files.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = $statusHandedBack
$wsZh.Range("B3").Value = $statusHandedBack

$wsZh.Range("E2").Value = "a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/8c0cc3595496c7b27b24961d96bb532b9af29e2d/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
Style-AsHyperlink $wsZh.Range("E2")

$wsZh.Range("F2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9e4be9cf36b2b0e4810178c6d2343355a701730a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
Style-AsHyperlink $wsZh.Range("F2")

$wsZh.Range("G2").Value = "2016-03-09 07:53:05"

$wsZh.Range("E3").Value = "a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/8c0cc3595496c7b27b24961d96bb532b9af29e2d/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
Style-AsHyperlink $wsZh.Range("E3")

$wsZh.Range("F3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9e4be9cf36b2b0e4810178c6d2343355a701730a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
Style-AsHyperlink $wsZh.Range("F3")

$wsZh.Range("G3").Value = "2016-03-09 07:53:05"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $statusHandedBack
$wsDe.Range("B3").Value = $statusHandedBack

$wsDe.Range("E2").Value = "a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/8c0cc3595496c7b27b24961d96bb532b9af29e2d/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
Style-AsHyperlink $wsDe.Range("E2")

$wsDe.Range("F2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a87bb3125443298b844ca7746b1a42d885e0d3f8/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
Style-AsHyperlink $wsDe.Range("F2")

$wsDe.Range("G2").Value = "2016-03-09 07:53:24"

$wsDe.Range("E3").Value = "a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/8c0cc3595496c7b27b24961d96bb532b9af29e2d/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
Style-AsHyperlink $wsDe.Range("E3")

$wsDe.Range("F3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a87bb3125443298b844ca7746b1a42d885e0d3f8/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
Style-AsHyperlink $wsDe.Range("F3")

$wsDe.Range("G3").Value = "2016-03-09 07:53:24"

Write-Host "Handback report generated."
